# "Implemented Edit Class and Delete class"
#
# Updates the ClassDetailsForm test-data sheet:
#   - Row 2 (AllMandateFieldsValidData): edited Topic/Dates/StaffName, added ExpectedResult
#   - Row 3: the old "AllMandateFieldsInValidData" case is replaced by a new
#            "ValidateInvalidMandatefields" case (edit-class negative test)
#   - Rows 4-9: StaffName + SelectClassDates refreshed to the new batch values
#   - Row 9: gets an (empty) ExpectedResult cell like its neighbours
#   - Row 10: a brand-new "AllMandateFieldsValidDataForEdit" case is appended
#             (delete-class / edit-class happy path)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- shared values used across many rows ----
$newDates = "03/20/2025,03/21/2025"
$newStaff = "Saranya M"

function Set-Cell($addr, $value) {
    $c = $ws.Range($addr)
    $c.Value = $value
    return $c
}

function Style-Body($c) {
    # plain body text: Arial 10, not bold (matches style used across the sheet, e.g. H column)
    $c.Font.Name = "Arial"
    $c.Font.Size = 10
    $c.Font.Bold = $false
}

function Style-Staff($c) {
    # StaffName column uses the larger (11pt) font
    $c.Font.Name = "Arial"
    $c.Font.Size = 11
    $c.Font.Bold = $false
}

# ===================== Row 2 =====================
$c = Set-Cell "C2" "Calculus3333"
Style-Body $c
$c = Set-Cell "E2" $newDates
Style-Body $c
$c = Set-Cell "G2" $newStaff
Style-Staff $c
$c = Set-Cell "L2" "Success"
Style-Body $c

# ===================== Row 3 (replaced test case) =====================
$c = Set-Cell "A3" "ValidateInvalidMandatefields"
Style-Body $c
$c = Set-Cell "B3" "Python101"
Style-Body $c
$c = Set-Cell "E3" $newDates
Style-Body $c
$c = Set-Cell "G3" $newStaff
Style-Staff $c
$c = Set-Cell "H3" "Active"
Style-Body $c
$c = Set-Cell "L3" "Error"
Style-Body $c

# ===================== Row 4 =====================
$c = Set-Cell "E4" $newDates
Style-Body $c
$c = Set-Cell "G4" $newStaff
Style-Staff $c

# ===================== Row 5 =====================
$c = Set-Cell "E5" $newDates
Style-Body $c
$c = Set-Cell "G5" $newStaff
Style-Staff $c

# ===================== Row 6 =====================
$c = Set-Cell "E6" $newDates
Style-Body $c
$c = Set-Cell "G6" $newStaff
Style-Staff $c

# ===================== Row 7 =====================
$c = Set-Cell "G7" $newStaff
Style-Staff $c

# ===================== Row 8 =====================
$c = Set-Cell "E8" $newDates
Style-Body $c
$c = Set-Cell "G8" $newStaff
Style-Staff $c

# ===================== Row 9 =====================
$c = Set-Cell "E9" $newDates
Style-Body $c
$c = Set-Cell "G9" $newStaff
Style-Staff $c
# empty styled ExpectedResult cell to match its neighbours
$c = $ws.Range("L9")
Style-Body $c

# ===================== Row 10 (new row) =====================
$ws.Rows.Item(10).RowHeight = 15.75

$c = Set-Cell "A10" "AllMandateFieldsValidDataForEdit"
Style-Body $c
$c = Set-Cell "C10" "Java1212"
Style-Body $c
$c = Set-Cell "E10" $newDates
Style-Body $c
$c = Set-Cell "G10" $newStaff
Style-Staff $c
$c = Set-Cell "H10" "Active"
Style-Body $c
$c = Set-Cell "L10" "Success"
Style-Body $c

# ---- sheet level touch-ups ----
$ws.Columns.Item(1).ColumnWidth = 30
$ws.Columns.Item(12).ColumnWidth = 17.42578125

$ws.Range("C2").Select()
